$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (excluding the header row) by the "Mark" column (D)
# in ascending order, as a stable sort - matches Data > Sort in Excel.
$dataRange = $ws.Range("A2:D13")
$sortKey = $ws.Range("D2:D13")
$dataRange.Sort($sortKey, 1)
